# The data block B2:Z26 (a 25x25 grid of calorimetry values) needs to be
# transposed in place: new[row, col] = old[col, row]. Row 1 (header values)
# and column A (header values) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rng = $ws.Range("B2:Z26")
$vals = $rng.Value()

$rows = $vals.GetLength(0)
$cols = $vals.GetLength(1)

$trans = New-Object 'object[,]' $rows, $cols
for ($i = 1; $i -le $rows; $i++) {
    for ($j = 1; $j -le $cols; $j++) {
        $trans[$i - 1, $j - 1] = $vals[$j, $i]
    }
}

$rng.Value = $trans
